# Fruta / hortaliza, semanal
# Reorders the weekly price records for rows 2-10 (Vega Monumental Concepción - Damasco)
# so that each row carries the correct date/variety/price data, per the updated dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44545
$ws.Cells.Item(2, 11).Value = 'Castle Brite'
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 100
$ws.Cells.Item(2, 14).Value = 18000
$ws.Cells.Item(2, 15).Value = 19000
$ws.Cells.Item(2, 16).Value = 18500
$ws.Cells.Item(2, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 1233
$ws.Cells.Item(2, 20).Value = 15

# Row 3
$ws.Cells.Item(3, 4).Value = 44545
$ws.Cells.Item(3, 11).Value = 'Castle Brite'
$ws.Cells.Item(3, 12).Value = 'Segunda'
$ws.Cells.Item(3, 13).Value = 50
$ws.Cells.Item(3, 14).Value = 17000
$ws.Cells.Item(3, 15).Value = 17000
$ws.Cells.Item(3, 16).Value = 17000
$ws.Cells.Item(3, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(3, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 19).Value = 1133
$ws.Cells.Item(3, 20).Value = 15

# Row 4
$ws.Cells.Item(4, 4).Value = 44579
$ws.Cells.Item(4, 11).Value = 'Modesto'
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 180
$ws.Cells.Item(4, 14).Value = 13000
$ws.Cells.Item(4, 15).Value = 14000
$ws.Cells.Item(4, 16).Value = 13444
$ws.Cells.Item(4, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(4, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 19).Value = 747
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value = 44187
$ws.Cells.Item(5, 11).Value = 'Dina'
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 100
$ws.Cells.Item(5, 14).Value = 15000
$ws.Cells.Item(5, 15).Value = 16000
$ws.Cells.Item(5, 16).Value = 15500
$ws.Cells.Item(5, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(5, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(5, 19).Value = 861
$ws.Cells.Item(5, 20).Value = 18

# Row 6
$ws.Cells.Item(6, 4).Value = 44159
$ws.Cells.Item(6, 11).Value = 'Castle Brite'
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 100
$ws.Cells.Item(6, 14).Value = 14000
$ws.Cells.Item(6, 15).Value = 15000
$ws.Cells.Item(6, 16).Value = 14500
$ws.Cells.Item(6, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(6, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 19).Value = 967
$ws.Cells.Item(6, 20).Value = 15

# Row 7
$ws.Cells.Item(7, 4).Value = 44189
$ws.Cells.Item(7, 11).Value = 'Dina'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 15000
$ws.Cells.Item(7, 15).Value = 16000
$ws.Cells.Item(7, 16).Value = 15500
$ws.Cells.Item(7, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(7, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(7, 19).Value = 1033
$ws.Cells.Item(7, 20).Value = 15

# Row 8
$ws.Cells.Item(8, 4).Value = 44189
$ws.Cells.Item(8, 11).Value = 'Dina'
$ws.Cells.Item(8, 12).Value = 'Segunda'
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 14).Value = 14000
$ws.Cells.Item(8, 15).Value = 14000
$ws.Cells.Item(8, 16).Value = 14000
$ws.Cells.Item(8, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 933
$ws.Cells.Item(8, 20).Value = 15

# Row 9
$ws.Cells.Item(9, 4).Value = 44559
$ws.Cells.Item(9, 11).Value = 'Modesto'
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 100
$ws.Cells.Item(9, 14).Value = 19000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 19500
$ws.Cells.Item(9, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(9, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(9, 19).Value = 1083
$ws.Cells.Item(9, 20).Value = 18

# Row 10
$ws.Cells.Item(10, 4).Value = 44559
$ws.Cells.Item(10, 11).Value = 'Modesto'
$ws.Cells.Item(10, 12).Value = 'Segunda'
$ws.Cells.Item(10, 13).Value = 50
$ws.Cells.Item(10, 14).Value = 18000
$ws.Cells.Item(10, 15).Value = 18000
$ws.Cells.Item(10, 16).Value = 18000
$ws.Cells.Item(10, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10, 19).Value = 1000
$ws.Cells.Item(10, 20).Value = 18
